# Update RF values in column I (rows 16-33) of Sheet1 from 4.3432 to 2.08
# per "Update of 2025 data and RF changes"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I16:I33").Value = 2.08
